# Split ISIC 05T06 into coal mining (ISIC 05) and oil and gas extraction (ISIC 06)
#
# The "SoCaOMSbRIC" sheet has one row of ISIC-code column headers (row 1)
# with a data row below (row 2). The column that used to hold the combined
# "ISIC 05T06" bucket needs to become two columns: "ISIC 05" and "ISIC 06".
# That combined bucket lived in column C, so we insert a new blank column at
# D (pushing the old D onward one column to the right) and then label the
# two resulting columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCaOMSbRIC")

# Insert a new column before D -- this shifts every column from D onward
# one column to the right, carrying over values/styles/formulas untouched.
$ws.Columns("D").Insert()

# Re-label the split columns and seed the new data cell.
$ws.Range("C1").Value = "ISIC 05"
$ws.Range("D1").Value = "ISIC 06"

# The original "ISIC 05T06" spending share was 0, split across both new
# buckets (C already holds 0 from the original column; D is the new blank
# cell created by the insert).
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
